$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "5, 2, 3, 8, 9, 3, 0, 5,"
$ws.Range("A2").Value = " 1, 1, 1, 0, 1, 1, 0, 1,"
$ws.Range("A3").Value = " 0, 0, 0, 0, 0, 2, 0, 1,"
$ws.Range("A4").Value = " 0, 0, 0, 0, 0, 0, 0, 0,"
$ws.Range("A5").Value = " 0, 0, 0, 1,-1, 0, 0, 0,"
$ws.Range("A6").Value = " 0, 0, 0, 0, 0, 0, 0, 0,"
$ws.Range("A7").Value = "-1,-1,-1,-1, 0,-1,-1,-1,"
$ws.Range("A8").Value = "-5,-2,-3, 0,-9,-3,-2,-5"

$excel.CalculateFull()
